$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 3894
$ws.Range("E2").Value = 26
$ws.Range("F2").Value = 57
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 10
$ws.Range("K2").Value = 1787
$ws.Range("L2").Value = 714
$ws.Range("M2").Value = 1073
$ws.Range("N2").Value = 1073
$ws.Range("P2").Value = 246
$ws.Range("Q2").Value = 126
$ws.Range("R2").Value = 61
$ws.Range("S2").Value = -109
$ws.Range("T2").Value = 45
$ws.Range("U2").Value = 81
$ws.Range("V2").Value = 3
$ws.Range("W2").Value = 0.67
$ws.Range("X2").Value = 0.26
$ws.Range("Y2").Value = 0.96
$ws.Range("Z2").Value = 0.5600000000000001
$ws.Range("AA2").Value = 66.52
$ws.Range("AB2").Value = 291.57
$ws.Range("AC2").Value = 17
$ws.Range("AD2").Value = 99.79000000000001
$ws.Range("AE2").Value = 1762
$ws.Range("AF2").Value = 0.96
$ws.Range("AG2").Value = 16
$ws.Range("AH2").Value = 0.9399999999999999
$ws.Range("AI2").Value = 93.43000000000001
$ws.Range("AJ2").Value = 60904132
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 3749
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 6
$ws.Range("G3").Value = 11
$ws.Range("H3").Value = 562
$ws.Range("I3").Value = 562
$ws.Range("K3").Value = 902
$ws.Range("L3").Value = 498
$ws.Range("M3").Value = 404
$ws.Range("N3").Value = 404
$ws.Range("P3").Value = 88
$ws.Range("Q3").Value = -13
$ws.Range("R3").Value = -152
$ws.Range("S3").Value = -14
$ws.Range("T3").Value = 13
$ws.Range("U3").Value = -26
$ws.Range("W3").Value = 0.17
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 76.12
$ws.Range("Z3").Value = 41.83
$ws.Range("AA3").Value = 123.11
$ws.Range("AB3").Value = 1152.59
$ws.Range("AC3").Value = 1380
$ws.Range("AD3").Value = 2.11
$ws.Range("AE3").Value = 1934
$ws.Range("AF3").Value = 1.51
$ws.Range("AG3").Value = 16
$ws.Range("AH3").Value = 0.54
$ws.Range("AI3").Value = 0.59
$ws.Range("AJ3").Value = 20905648
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("V3").ClearContents()

# Row 4
$ws.Range("D4").Value = 3847
$ws.Range("E4").Value = 51
$ws.Range("F4").Value = 51
$ws.Range("G4").Value = 49
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 36
$ws.Range("K4").Value = 979
$ws.Range("L4").Value = 539
$ws.Range("M4").Value = 440
$ws.Range("N4").Value = 440
$ws.Range("P4").Value = 88
$ws.Range("Q4").Value = 37
$ws.Range("R4").Value = -30
$ws.Range("S4").Value = -3
$ws.Range("T4").Value = 42
$ws.Range("U4").Value = -5
$ws.Range("W4").Value = 1.31
$ws.Range("X4").Value = 0.95
$ws.Range("Y4").Value = 8.609999999999999
$ws.Range("Z4").Value = 3.86
$ws.Range("AA4").Value = 122.33
$ws.Range("AB4").Value = 1191.57
$ws.Range("AC4").Value = 174
$ws.Range("AD4").Value = 13.54
$ws.Range("AE4").Value = 2107
$ws.Range("AF4").Value = 1.12
$ws.Range("AG4").Value = 16
$ws.Range("AH4").Value = 0.67
$ws.Range("AI4").Value = 9.06
$ws.Range("AJ4").Value = 20905648
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
$ws.Range("V4").ClearContents()

# Row 5
$ws.Range("D5").Value = 3553
$ws.Range("E5").Value = 23
$ws.Range("F5").Value = 23
$ws.Range("G5").Value = 39
$ws.Range("H5").Value = 29
$ws.Range("I5").Value = 29
$ws.Range("K5").Value = 919
$ws.Range("L5").Value = 447
$ws.Range("M5").Value = 472
$ws.Range("N5").Value = 472
$ws.Range("P5").Value = 88
$ws.Range("Q5").Value = 16
$ws.Range("R5").Value = -13
$ws.Range("S5").Value = -3
$ws.Range("T5").Value = 16
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0.65
$ws.Range("X5").Value = 0.8100000000000001
$ws.Range("Y5").Value = 6.3
$ws.Range("Z5").Value = 3.03
$ws.Range("AA5").Value = 94.66
$ws.Range("AB5").Value = 1228.64
$ws.Range("AC5").Value = 137
$ws.Range("AD5").Value = 12.66
$ws.Range("AE5").Value = 2258
$ws.Range("AF5").Value = 0.77
$ws.Range("AG5").Value = 24
$ws.Range("AH5").Value = 1.36
$ws.Range("AI5").Value = 17.21
$ws.Range("AJ5").Value = 20905648
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 4342
$ws.Range("E6").Value = 39
$ws.Range("F6").Value = 39
$ws.Range("G6").Value = 23
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = 15
$ws.Range("K6").Value = 1268
$ws.Range("L6").Value = 705
$ws.Range("M6").Value = 563
$ws.Range("N6").Value = 563
$ws.Range("P6").Value = 143
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = -114
$ws.Range("S6").Value = 201
$ws.Range("T6").Value = 42
$ws.Range("U6").Value = -42
$ws.Range("V6").Value = 112
$ws.Range("W6").Value = 0.9
$ws.Range("X6").Value = 0.34
$ws.Range("Y6").Value = 2.81
$ws.Range("Z6").Value = 1.33
$ws.Range("AA6").Value = 125.24
$ws.Range("AB6").Value = 788.1900000000001
$ws.Range("AC6").Value = 67
$ws.Range("AD6").Value = 19.75
$ws.Range("AE6").Value = 2049
$ws.Range("AF6").Value = 0.65
$ws.Range("AG6").Value = 30
$ws.Range("AH6").Value = 2.26
$ws.Range("AI6").Value = 56.67
$ws.Range("AJ6").Value = 27477425

# Rows 7-9: clear all data columns except A, B, C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
